# The edit swaps the Id / Antal / Ost / Nord values between row 2 and row 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 values (previously held by row 4)
$ws.Range("A2").Value2 = 111651416
$ws.Range("I2").Value2 = 10
$ws.Range("Q2").Value2 = 573987.3009507703
$ws.Range("R2").Value2 = 6403998.840001023

# New row 4 values (previously held by row 2)
$ws.Range("A4").Value2 = 111651333
$ws.Range("I4").Value2 = 25
$ws.Range("Q4").Value2 = 574013.3982996774
$ws.Range("R4").Value2 = 6403974.780079928
